$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values: Año / Evento / Fecha for the 2017 and 2021 elections ---
$ws.Range("A14").Value = 2017
$ws.Range("B14").Value = "PASO"
$ws.Range("C14").Value = 42960

$ws.Range("B15").Value = "Generales"
$ws.Range("C15").Value = 43030

$ws.Range("A16").Value = 2021
$ws.Range("B16").Value = "PASO"
$ws.Range("C16").Value = 44451

$ws.Range("B17").Value = "Generales"
$ws.Range("C17").Value = 44514

# --- Formats: reuse the existing "Evento" cell style (font/border/wrap) for column B ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B14:B17").PasteSpecial(-4122) | Out-Null

# --- Formats: apply the date number format to column C (new plain date style) ---
$ws.Range("C14").NumberFormat = "mm-dd-yy"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15:C17").PasteSpecial(-4122) | Out-Null

$ws.Range("B16:B17").Select() | Out-Null
